# Update cryptos list: refresh prices / volumes and fix Toncoin/RenderToken row order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.879.18"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "'2.929.79"
$ws.Range("E3").Value = "  +3.12%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'352.39"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'112.26"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Value = "'39.43"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "'0.0885"
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'20.14"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "'3.387.92"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "'2.934.15"
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "'0.985"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'51.930.73"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -4.43%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "  +6.29%  "
$ws.Range("D22").Value = "'0.0₃0984"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").Value = "'71.22"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'268.44"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +11.38%  "
$ws.Range("D27").Value = "'26.97"
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = "  +14.96%  "
$ws.Range("E30").Value = "  +16.28%  "
$ws.Range("D31").Value = "'10.58"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'37.21"
$ws.Range("E32").Value = "  -4.67%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'2.26"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'6.23"
$ws.Range("E34").Value = "  +9.88%  "
$ws.Range("D35").Value = "'52.91"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("D39").Value = "'18.73"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("E41").Value = "  +6.13%  "
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").Value = "'23.16"
$ws.Range("E43").Value = "  +4.22%  "
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("D46").Value = "'3.50"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "'2.174.58"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'111.73"
$ws.Range("E48").Value = "  -8.39%  "
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  +10.29%  "
$ws.Range("D51").Value = "'0.945"
$ws.Range("E51").Value = "  -3.19%  "
